$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F3").Value  = "30,46 TL - 60,94 TL - 609,43 TL"
$ws.Range("K3").Value  = "30,46 TL - 60,94 TL - 609,43 TL"

$ws.Range("F4").Value  = "30,46 TL - 60,94 TL - 609,43 TL"
$ws.Range("K4").Value  = "30,46 TL - 60,94 TL - 609,43 TL"

$ws.Range("F5").Value  = "30,46 TL - 60,94 TL - 609,43 TL"
$ws.Range("K5").Value  = "30,46 TL - 60,94 TL - 609,43 TL"

$ws.Range("K6").Value  = "6,09 TL - 12,19 TL - 152,35 TL"

$ws.Range("F8").Value  = "15,23 TL - 30,47 TL - 304,71 TL"
$ws.Range("K8").Value  = "15,23 TL - 30,47 TL - 304,71 TL"

$ws.Range("F9").Value  = "15,23 TL - 30,47 TL - 304,71 TL"
$ws.Range("K9").Value  = "15,23 TL - 30,47 TL - 304,71 TL"

$ws.Range("F10").Value = "15,23 TL - 30,47 TL - 304,71 TL"
$ws.Range("K10").Value = "15,23 TL - 30,47 TL - 304,71 TL"

$ws.Range("K11").Value = "3,05 TL - 6,09 TL - 76,17 TL"

$ws.Range("K12").Value = "WU: ,USD–; Diğer: 404,16 TL–3.403,42 TL"

$ws.Range("F13").Value = "Hesaba: Asgari 300 TL | Azami 3.080 TL"
$ws.Range("K13").Value = "Hesaba: Asgari 1 TL | Azami 865,75 TL"

$ws.Range("F14").Value = "1.952,38 TL - 9.523,81 TL"
$ws.Range("K14").Value = "914,14 TL - 4.265,98 TL"
